# Atualização responsabilidade incremento 3 [Douglas Giordano]
$wb = $excel.ActiveWorkbook

# --- Sheet "Requisitos": mark 4.8 as Finalizada, and move 5.2/5.3 to Incremento 3 under Douglas ---
$reqs = $wb.Worksheets.Item("Requisitos")
$reqs.Range("F24").Value = "Finalizada"
$reqs.Range("D29").Value = 3
$reqs.Range("E29").Value = "Douglas"
$reqs.Range("D30").Value = 3
$reqs.Range("E30").Value = "Douglas"

# --- Sheet "Incremento 2": mark 4.8 as Finalizada ---
$inc2 = $wb.Worksheets.Item("Incremento 2")
$inc2.Range("H11").Value = "Finalizada"

# --- Sheet "Incremento 3": add the two newly promoted use cases for Douglas ---
$inc3 = $wb.Worksheets.Item("Incremento 3")
$inc3.Range("A6").NumberFormat = "@"
$inc3.Range("A7").NumberFormat = "@"
$inc3.Range("A6").Value = "5.2"
$inc3.Range("B6").Value = "Criar ata de julgamento da prova de títulos"
$inc3.Range("C6").Value = 4
$inc3.Range("D6").Value = 2
$inc3.Range("E6").Value = "Douglas"
$inc3.Range("F6").Value = "Pendente"
$inc3.Range("G6").Value = "Pendente"
$inc3.Range("H6").Value = "Em andamento"
$inc3.Range("I6").Value = 4

$inc3.Range("A7").Value = "5.3"
$inc3.Range("B7").Value = "Criar recibo de devolução de documentação"
$inc3.Range("C7").Value = 4
$inc3.Range("D7").Value = 2
$inc3.Range("E7").Value = "Douglas"
$inc3.Range("F7").Value = "Pendente"
$inc3.Range("G7").Value = "Pendente"
$inc3.Range("H7").Value = "Pendente"
$inc3.Range("I7").Value = 3

# Resource-hours summary: Douglas planned hours for this increment
$inc3.Range("M4").Value = 7

# Status summary for the "Casos de Uso Atrasados" block
$inc3.Range("M16").Value = 0
$inc3.Range("M17").Value = 0
$inc3.Range("M18").Value = 3
$inc3.Range("M19").Value = 1
